$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at 65 (pushes old rows 65-80 down to 66-81) ---
$ws.Rows("65:65").Insert()

# --- Insert two new rows before the last row (old row 80, now row 81) ---
$ws.Rows("81:82").Insert()

# Set values in the order that matches the shared-string insertion order
# expected by the target workbook (sharedStrings.xml appends in first-use
# order): "Show annual summary..." (115), then "...so that..." (116),
# then "...Flat Plate..." (117).
$ws.Range("A82").Value = "Future"
$ws.Range("B82").Value = "Show annual summary statistics (irradiance, wind speed, ambient temp) on the resource page"

$ws.Range("A81").Value = "Future"
$ws.Range("B81").Value = "Fix output variable group names in SSC so that they show up pretty in outputs browser"

$ws.Range("A65").Value = "Not done"
$ws.Range("B65").Value = "Fix output variable group names in SSC: Flat Plate, PVWatts, Wind, CSP Trough & Towers, Res, Comm"

# --- Update view/selection to match target ---
$ws.Range("A66").Select()
$excel.ActiveWindow.ScrollRow = 59
$excel.ActiveWindow.ScrollColumn = 1
